$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.192.97'
$ws.Range("E2").Value = '  +2.56%  '

$ws.Range("D3").Value = '2.086.36'
$ws.Range("E3").Value = '  +3.56%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.38'
$ws.Range("E5").Value = '  +2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.661'

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.54'
$ws.Range("E8").Value = '  +21.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.60'
$ws.Range("E9").Value = '  +1.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.373'
$ws.Range("E10").Value = '  +4.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0741'
$ws.Range("E11").Value = '  +4.22%  '

$ws.Range("E12").Value = '  +8.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.01'
$ws.Range("E13").Value = '  +5.90%  '

$ws.Range("D14").Value = '2.394.07'
$ws.Range("E14").Value = '  +3.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.829'
$ws.Range("E15").Value = '  +3.74%  '

$ws.Range("D16").Value = '2.091.75'
$ws.Range("E16").Value = '  +3.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.14'

$ws.Range("D18").Value = '37.163.05'
$ws.Range("E18").Value = '  +2.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.62'
$ws.Range("E19").Value = '  +2.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.55'
$ws.Range("E20").Value = '  +15.04%  '

$ws.Range("D21").Value = '0.0₃0843'
$ws.Range("E21").Value = '  +4.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.18'
$ws.Range("E22").Value = '  +1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.18'
$ws.Range("E23").Value = '  +6.60%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("E25").Value = '  +2.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.13'
$ws.Range("E26").Value = '  +3.21%  '

$ws.Range("E27").Value = '  +7.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.63'
$ws.Range("E28").Value = '  +5.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("E29").Value = '  +3.09%  '

$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.18'
$ws.Range("E31").Value = '  +8.16%  '

$ws.Range("E32").Value = '  +27.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("E33").Value = '  +4.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0615'
$ws.Range("E34").Value = '  +6.78%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0911'
$ws.Range("E35").Value = '  +7.38%  '

$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.13'
$ws.Range("E37").Value = '  +4.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.83'
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("E39").Value = '  +4.71%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.34'
$ws.Range("E40").Value = '  +2.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.16'
$ws.Range("E41").Value = '  +16.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("E42").Value = '  +5.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.16'
$ws.Range("E43").Value = '  +6.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.04'
$ws.Range("E44").Value = '  +3.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0941'
$ws.Range("E45").Value = '  +15.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.33'
$ws.Range("E46").Value = '  +122.24%  '

$ws.Range("E47").Value = '  +0.69%  '

$ws.Range("D48").Value = '1.319.00'
$ws.Range("E48").Value = '  +1.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.93'
$ws.Range("E49").Value = '  +5.63%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.99'
$ws.Range("E50").Value = '  +15.37%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("E51").Value = '  +6.23%  '
